$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "wong3"

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("B5").Value = 85
$ws.Range("C5").Value = 85
$ws.Range("E5").Value = 95
$ws.Range("F5").Value = 95
$ws.Range("H5").Value = 109
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 12
$ws.Range("F6").Value = 12
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 7
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 4
$ws.Range("B12").Value = 12
$ws.Range("C12").Value = 12
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 7
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("B16").Value = 12
$ws.Range("C16").Value = 12
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 19
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 4
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 5
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 4
$ws.Range("B23").Value = 57
$ws.Range("C23").Value = 56
$ws.Range("E23").Value = 73
$ws.Range("F23").Value = 73
$ws.Range("H23").Value = 78
$ws.Range("E25").Value = 6
$ws.Range("F25").Value = 6
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 2
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 46
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = 2
$ws.Range("H31").Value = 15
$ws.Range("B34").Value = 2
$ws.Range("C34").Value = 2
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = 3
$ws.Range("B39").Value = 5
$ws.Range("C39").Value = 5
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 7
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 1
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 1
$ws.Range("B44").Value = 8
$ws.Range("C44").Value = 8
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 15
$ws.Range("B45").Value = 3
$ws.Range("C45").Value = 3
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 4
$ws.Range("B47").Value = 1
$ws.Range("C47").Value = 1
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = 1
$ws.Range("B48").Value = 1
$ws.Range("C48").Value = 1
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 1
$ws.Range("F51").Value = 13
$ws.Range("H51").Value = 21
$ws.Range("B53").Value = 3
$ws.Range("C53").Value = 2
$ws.Range("F53").Value = 6
$ws.Range("H53").Value = 19
$ws.Range("E55").Value = 7
$ws.Range("F55").Value = 7
$ws.Range("E57").Value = 2
$ws.Range("F57").Value = 2
$ws.Range("E59").Value = 1
$ws.Range("F59").Value = 1
$ws.Range("B62").Value = 2
$ws.Range("C62").Value = 2
$ws.Range("E62").Value = 2
$ws.Range("F62").Value = 2
$ws.Range("E63").Value = 1
$ws.Range("F63").Value = 1
$ws.Range("E66").Value = 2
$ws.Range("F66").Value = 2
$ws.Range("B67").Value = 2
$ws.Range("C67").Value = 2
$ws.Range("E67").Value = 7
$ws.Range("F67").Value = 7
$ws.Range("B68").Value = 17
$ws.Range("C68").Value = 17
$ws.Range("E68").Value = 30
$ws.Range("F68").Value = 30
$ws.Range("B69").Value = 1
$ws.Range("C69").Value = 1
$ws.Range("E69").Value = 1
$ws.Range("F69").Value = 1
$ws.Range("B70").Value = 1
$ws.Range("C70").Value = 1
$ws.Range("E71").Value = 2
$ws.Range("F71").Value = 2
$ws.Range("B75").Value = 4
$ws.Range("C75").Value = 4
$ws.Range("E75").Value = 7
$ws.Range("F75").Value = 7
$ws.Range("B76").Value = 4
$ws.Range("C76").Value = 4
$ws.Range("E76").Value = 10
$ws.Range("F76").Value = 10
$ws.Range("B79").Value = 6
$ws.Range("C79").Value = 5
$ws.Range("E79").Value = 13
$ws.Range("F79").Value = 9
$ws.Range("B82").Value = 2
$ws.Range("C82").Value = 2
$ws.Range("E82").Value = 8
$ws.Range("F82").Value = 8
$ws.Range("B85").Value = 2
$ws.Range("C85").Value = 2
$ws.Range("B86").Value = 6
$ws.Range("C86").Value = 6
$ws.Range("E86").Value = 8
$ws.Range("F86").Value = 8
$ws.Range("B89").Value = 4
$ws.Range("C89").Value = 4
$ws.Range("E89").Value = 14
$ws.Range("F89").Value = 14
$ws.Range("B90").Value = 1
$ws.Range("C90").Value = 1
$ws.Range("B91").Value = 22
$ws.Range("C91").Value = 22
$ws.Range("E91").Value = 39
$ws.Range("F91").Value = 39
$ws.Range("H91").Value = 49
$ws.Range("B92").Value = 2
$ws.Range("C92").Value = 2
$ws.Range("E92").Value = 2
$ws.Range("F92").Value = 2
$ws.Range("B93").Value = 2
$ws.Range("C93").Value = 2
$ws.Range("E93").Value = 3
$ws.Range("F93").Value = 3
